$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab
$ws.Name = "ShearF"

# Fix tiny precision update on H13
$ws.Range("H13").Value = 0.9940111371242357

# Add new row 16: HexGrid-60degTilt5degRes data (HKL index 14)
$ws.Range("A16").Value = 14
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9957675914803306
$ws.Range("D16").Value = 1.64344820480699
$ws.Range("E16").Value = 0.8079618251480029
$ws.Range("F16").Value = 0.9957675914803306
$ws.Range("G16").Value = 1.280637724328747
$ws.Range("H16").Value = 0.6227848949555937
$ws.Range("I16").Value = 0.8534831158739909
$ws.Range("J16").Value = 1.64344820480699
$ws.Range("K16").Value = 1.225705014977496
$ws.Range("L16").Value = 1.110736303228913
$ws.Range("M16").Value = 1.034013892765609
